# The workbook was previously produced by an ad-hoc Excel writer that had
# left explicit column-width overrides (cols) on the sheet. The commit
# switches the writer to a pandas DataFrame based export, which re-writes
# the worksheet from scratch with no manual column widths, and appends a
# second data row (a dummy placeholder row so the PK calc always has a
# value to work with).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the manual column-width formatting (min=2,4,5,6,8,9) that the old
# writer had set, so the sheet falls back to default column widths - same
# as a freshly written pandas/openpyxl sheet would have.
$ws.Range("A1:Z1").EntireColumn.Delete()

# Re-write the header row (row 1) exactly as before.
$ws.Range("A1").Value = "PK"
$ws.Range("B1").Value = "CARD_NAME"
$ws.Range("C1").Value = "SERIES"
$ws.Range("D1").Value = "SET_NAME"
$ws.Range("E1").Value = "NUM_IN_SET"
$ws.Range("F1").Value = "SET_TOTAL"
$ws.Range("G1").Value = "FOIL"
$ws.Range("H1").Value = "CONDITION"
$ws.Range("I1").Value = "EBAY_TITLE"

# Add the new placeholder data row (row 2).
$ws.Range("A2").Value = 0
$ws.Range("B2:H2").Value = "t"
$ws.Range("I2").Value = "There needs to be a value in this PK for the pk_calc to work"

# Match the selection left behind by the editing session.
$ws.Range("R7").Select()
